$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Stemmed"
$ws.Range("G2").Value = 10000
$ws.Range("G3").Value = "No Null"
$ws.Range("G4").Value = 0.36879400000000001
$ws.Range("G5").Value = 0.59055100000000005
$ws.Range("G6").Value = 0.69230800000000003
$ws.Range("G7").Value = 0.28571400000000002
$ws.Range("G8").Value = 0.68143500000000001
$ws.Range("G9").Value = 0.760355
$ws.Range("G10").Value = 0.30713499999999999
$ws.Range("G11").Value = 0.61049299999999995
$ws.Range("G12").Value = 0.84615399999999996
$ws.Range("G13").Value = 0.35745300000000002
$ws.Range("G14").Value = 0.67705400000000004
$ws.Range("G15").Value = 0.60650899999999996
$ws.Range("G16").Value = 0.35452299999999998
$ws.Range("G17").Value = 0.61458299999999999
$ws.Range("G18").Value = 0.68934899999999999

$ws.Range("G2").NumberFormat = $ws.Range("F2").NumberFormat

$ws.Range("G18").Select()
